$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, [string]$addr, [double]$val)
    $ws.Range($addr).Value = $val
}

function Clear-CellValue {
    param($ws, [string]$addr)
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H17" 889.2222
Set-CellValue $ws "J17" 889.2222
Set-CellValue $ws "L17" 2667.6666
Set-CellValue $ws "N17" -3003.6666
Set-CellValue $ws "H28" 10079.762
Set-CellValue $ws "I28" 3985.5
Set-CellValue $ws "J28" 22268.285
Set-CellValue $ws "K28" 3985.5
Set-CellValue $ws "L28" 22268.285
Set-CellValue $ws "M28" -3500.5
Set-CellValue $ws "N28" -23238.285
Set-CellValue $ws "H106" 2685.5386
Set-CellValue $ws "I106" 3664.5715
Set-CellValue $ws "K106" 3664.5715
Set-CellValue $ws "M106" -3033.5715
Set-CellValue $ws "H116" 6966.0
Set-CellValue $ws "I116" 7049.1665
Set-CellValue $ws "K116" 7049.1665
Set-CellValue $ws "M116" -3607.1665
Set-CellValue $ws "H137" 6630.375
Set-CellValue $ws "I137" 3433.0
Set-CellValue $ws "J137" 8548.8
Set-CellValue $ws "K137" 10299.0
Set-CellValue $ws "L137" 25646.4
Set-CellValue $ws "M137" -7749.0
Set-CellValue $ws "N137" -30746.4
Set-CellValue $ws "H138" 2626.475
Set-CellValue $ws "I138" 2193.5625
Set-CellValue $ws "J138" 2915.0833
Set-CellValue $ws "K138" 6580.6875
Set-CellValue $ws "L138" 8745.249899999999
Set-CellValue $ws "M138" -1440.6875
Set-CellValue $ws "N138" -19025.2499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H32" 3994.4575
Set-CellValue $ws "I32" 1581.48
Set-CellValue $ws "K32" 1581.48
Set-CellValue $ws "M32" -1294.48
Set-CellValue $ws "H122" 2333.0
Set-CellValue $ws "I122" 2333.0
Set-CellValue $ws "K122" 6999.0
Set-CellValue $ws "M122" -4549.0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H58" 88420.25
Set-CellValue $ws "J58" 92482.29
Set-CellValue $ws "L58" 92482.29
Set-CellValue $ws "N58" -93070.29
Set-CellValue $ws "H105" 3249.3845
Set-CellValue $ws "I105" 3224.3
Set-CellValue $ws "K105" 3224.3
Set-CellValue $ws "M105" -1477.3
Set-CellValue $ws "H134" 2465.9023
Set-CellValue $ws "I134" 1920.2413
Set-CellValue $ws "K134" 5760.7239
Set-CellValue $ws "M134" -3225.7239

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H58" 3599.2144
Set-CellValue $ws "I58" 3599.2144
Set-CellValue $ws "K58" 3599.2144
Set-CellValue $ws "M58" -3396.2144
Set-CellValue $ws "H62" 13499.5
Set-CellValue $ws "I62" 13499.5
Set-CellValue $ws "K62" 13499.5
Set-CellValue $ws "M62" -12875.5
Set-CellValue $ws "H65" 13499.5
Set-CellValue $ws "I65" 13499.5
Set-CellValue $ws "K65" 67497.5
Set-CellValue $ws "M65" -64377.5
Set-CellValue $ws "H99" 24280.0
Set-CellValue $ws "I99" 6334.6665
Set-CellValue $ws "K99" 6334.6665
Set-CellValue $ws "M99" -4836.6665
Set-CellValue $ws "H126" 24280.0
Set-CellValue $ws "I126" 6334.6665
Set-CellValue $ws "K126" 19003.9995
Set-CellValue $ws "M126" -16533.9995
Set-CellValue $ws "H134" 2117.7144
Set-CellValue $ws "I134" 1905.0
Set-CellValue $ws "K134" 5715.0
Set-CellValue $ws "M134" -3180.0
Set-CellValue $ws "H136" 3599.2144
Set-CellValue $ws "I136" 3599.2144
Set-CellValue $ws "K136" 10797.6432
Set-CellValue $ws "M136" -8247.643199999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H22" 743.0
Set-CellValue $ws "I22" 800.0
Set-CellValue $ws "J22" 714.5
Set-CellValue $ws "K22" 2400.0
Set-CellValue $ws "L22" 2143.5
Set-CellValue $ws "M22" -2231.0
Set-CellValue $ws "N22" -2481.5
Set-CellValue $ws "H27" 743.0
Set-CellValue $ws "I27" 800.0
Set-CellValue $ws "J27" 714.5
Set-CellValue $ws "K27" 2400.0
Set-CellValue $ws "L27" 2143.5
Set-CellValue $ws "M27" -2298.0
Set-CellValue $ws "N27" -2347.5
Set-CellValue $ws "H40" 882.5263
Set-CellValue $ws "I40" 10.5
Set-CellValue $ws "J40" 5533.3335
Set-CellValue $ws "K40" 42.0
Set-CellValue $ws "L40" 22133.334
Set-CellValue $ws "M40" 27.0
Set-CellValue $ws "N40" -22271.334
Set-CellValue $ws "H51" 2081.25
Set-CellValue $ws "I51" 2064.6
Set-CellValue $ws "J51" 2109.0
Set-CellValue $ws "K51" 6193.799999999999
Set-CellValue $ws "L51" 6327.0
Set-CellValue $ws "M51" -5733.799999999999
Set-CellValue $ws "N51" -7247.0
Set-CellValue $ws "H54" 2333.3333
Set-CellValue $ws "H80" 14697.6875
Set-CellValue $ws "I80" 20324.375
Set-CellValue $ws "J80" 9071.0
Set-CellValue $ws "K80" 60973.125
Set-CellValue $ws "L80" 27213.0
Set-CellValue $ws "M80" -60037.125
Set-CellValue $ws "N80" -29085.0
Set-CellValue $ws "H83" 14697.6875
Set-CellValue $ws "I83" 20324.375
Set-CellValue $ws "J83" 9071.0
Set-CellValue $ws "K83" 182919.375
Set-CellValue $ws "L83" 81639.0
Set-CellValue $ws "M83" -178239.375
Set-CellValue $ws "N83" -90999.0
Set-CellValue $ws "H131" 223481.2
Set-CellValue $ws "I131" 556318.0
Set-CellValue $ws "K131" 1668954.0
Set-CellValue $ws "M131" -1663914.0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H64" 63423.668
Set-CellValue $ws "H67" 63423.668
Set-CellValue $ws "H102" 3834.0625
Set-CellValue $ws "I102" 3889.9167
Set-CellValue $ws "K102" 3889.9167
Set-CellValue $ws "M102" -2267.9167
Set-CellValue $ws "H122" 2511.2
Set-CellValue $ws "I122" 2364.0667
Set-CellValue $ws "J122" 2658.3333
Set-CellValue $ws "K122" 7092.2001
Set-CellValue $ws "L122" 7974.999899999999
Set-CellValue $ws "M122" -4642.2001
Set-CellValue $ws "N122" -12874.9999
Set-CellValue $ws "H132" 5502.75
Set-CellValue $ws "J132" 0.0
Set-CellValue $ws "L132" 0.0
Clear-CellValue $ws "N132"

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws "H7" 3627.5833
Set-CellValue $ws "I7" 3627.5833
Set-CellValue $ws "J7" 0.0
Set-CellValue $ws "K7" 3627.5833
Set-CellValue $ws "L7" 0.0
Set-CellValue $ws "M7" -3515.5833
Clear-CellValue $ws "N7"
Set-CellValue $ws "H9" 396.5
Set-CellValue $ws "I9" 622.6667
Set-CellValue $ws "J9" 170.33333
Set-CellValue $ws "K9" 622.6667
Set-CellValue $ws "L9" 170.33333
Set-CellValue $ws "M9" -398.6667
Set-CellValue $ws "N9" -618.3333299999999
Set-CellValue $ws "H10" 5500.0
Set-CellValue $ws "I10" 10000.0
Set-CellValue $ws "J10" 1000.0
Set-CellValue $ws "K10" 10000.0
Set-CellValue $ws "L10" 1000.0
Set-CellValue $ws "M10" -9860.0
Set-CellValue $ws "N10" -1280.0
Set-CellValue $ws "H13" 2500.0
Set-CellValue $ws "I13" 2500.0
Set-CellValue $ws "K13" 2500.0
Set-CellValue $ws "M13" -2360.0
Set-CellValue $ws "H40" 3972.8333
Set-CellValue $ws "I40" 3176.5625
Set-CellValue $ws "J40" 4882.857
Set-CellValue $ws "K40" 3176.5625
Set-CellValue $ws "L40" 4882.857
Set-CellValue $ws "M40" -3040.5625
Set-CellValue $ws "N40" -5154.857
Set-CellValue $ws "H46" 1474.1111
Set-CellValue $ws "I46" 1317.0
Set-CellValue $ws "K46" 1317.0
Set-CellValue $ws "M46" -1129.0
Set-CellValue $ws "H74" 39154.0
Set-CellValue $ws "J74" 39984.8
Set-CellValue $ws "L74" 39984.8
Set-CellValue $ws "N74" -41980.8
Set-CellValue $ws "H77" 39154.0
Set-CellValue $ws "J77" 39984.8
Set-CellValue $ws "L77" 119954.4
Set-CellValue $ws "N77" -129938.4
Set-CellValue $ws "H126" 3627.5833
Set-CellValue $ws "I126" 3627.5833
Set-CellValue $ws "J126" 0.0
Set-CellValue $ws "K126" 10882.7499
Set-CellValue $ws "L126" 0.0
Set-CellValue $ws "M126" -8412.749899999999
Clear-CellValue $ws "N126"
Set-CellValue $ws "H132" 3507.25
Set-CellValue $ws "I132" 3409.6843
Set-CellValue $ws "J132" 3713.2222
Set-CellValue $ws "K132" 10229.0529
Set-CellValue $ws "L132" 11139.6666
Set-CellValue $ws "M132" -7699.052899999999
Set-CellValue $ws "N132" -16199.6666
Set-CellValue $ws "H136" 3683.25
Set-CellValue $ws "I136" 3800.0
Set-CellValue $ws "J136" 3333.0
Set-CellValue $ws "K136" 11400.0
Set-CellValue $ws "L136" 9999.0
Set-CellValue $ws "M136" -8850.0
Set-CellValue $ws "N136" -15099.0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H122" 3336.5
Set-CellValue $ws "J122" 3512.0833
Set-CellValue $ws "L122" 10536.2499
Set-CellValue $ws "N122" -15436.2499
Set-CellValue $ws "H126" 5627.722
Set-CellValue $ws "I126" 5237.4165
Set-CellValue $ws "J126" 6408.3335
Set-CellValue $ws "K126" 15712.2495
Set-CellValue $ws "L126" 19225.0005
Set-CellValue $ws "M126" -13242.2495
Set-CellValue $ws "N126" -24165.0005
Set-CellValue $ws "H132" 4092.8057
Set-CellValue $ws "I132" 4090.0293
Set-CellValue $ws "K132" 12270.0879
Set-CellValue $ws "M132" -9740.0879
Set-CellValue $ws "H136" 4999.2
Set-CellValue $ws "I136" 4999.2
Set-CellValue $ws "K136" 14997.6
Set-CellValue $ws "M136" -12447.6
